$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.890.83"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "'3.691.98"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'647.65"
$ws.Range("E5").Value = "  -4.25%  "

$ws.Range("D6").Value = "'161.94"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +1.07%  "

$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").Value = "'7.21"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").Value = "'0.0000233"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "'4.315.63"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "'32.77"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").Value = "'3.679.57"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "'69.883.45"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").Value = "'6.52"

$ws.Range("D20").Value = "'10.40"
$ws.Range("E20").Value = "  +6.05%  "

$ws.Range("D21").Value = "'471.99"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").Value = "'80.13"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").Value = "'3.839.01"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  +1.10%  "

$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").Value = "'9.15"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "  -1.95%  "

$ws.Range("D31").Value = "'2.01"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").Value = "'0.169"
$ws.Range("E32").Value = "  +3.83%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("D34").Value = "'6.53"

$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").Value = "'3.688.87"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D39").Value = "'180.91"
$ws.Range("E39").Value = "  +8.46%  "

$ws.Range("E40").Value = "  -5.05%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.23"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").Value = "'0.0904"
$ws.Range("E43").Value = "  +0.30%  "

$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").Value = "'2.87"
$ws.Range("E45").Value = "  +3.59%  "

$ws.Range("D46").Value = "'29.44"
$ws.Range("E46").Value = "  +4.67%  "

$ws.Range("D47").Value = "'46.71"
$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("D48").Value = "'0.000274"
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("E49").Value = "  -3.06%  "

$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("E51").Value = "  -3.33%  "
